# "cierre 30 Dic 22" - close out week 52 (26 Dic - 01 Ene) payroll receipt.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("recibos")
$ws.Activate()

# Week label (B9 is the only literal cell; H9/B28/H28/B46 all reference it
# via formulas, so they update automatically on recalculation).
$ws.Range("B9").Value = "SEMANA  52   DEL    26      Al   01   DE    ENERO          2022"

# Updated payroll figures for the new week.
$ws.Range("K22").Value = 840
$ws.Range("E23").Value = 2300
$ws.Range("E43").Value = 0

# Restore the view state: scrolled up a bit, and the active cell moved.
$ws.Range("I40").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
